$d = $word.ActiveDocument

# 1. Ref-AB1CD2 -> Ref-f695660
$d.Content.Find.Execute("Ref-AB1CD2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-f695660", 2)

# 2. (Ref-A1B2C3) -> (Smith)
$d.Content.Find.Execute("(Ref-A1B2C3)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Smith)", 2)

# 3. Ref-J7X2BZ -> Ref-f155364
$d.Content.Find.Execute("Ref-J7X2BZ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-f155364", 2)

# 4. Ref-AB12CD -> Ref-s015350
$d.Content.Find.Execute("Ref-AB12CD", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s015350", 2)

# 5. Ref-EF34GH -> Ref-s015350
$d.Content.Find.Execute("Ref-EF34GH", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-s015350", 2)

# 6. Ref-DJ74KL -> Ref-f977484
$d.Content.Find.Execute("Ref-DJ74KL", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ref-f977484", 2)
